# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (column E, rows 16-22) is reordered: the
# database of periods is refreshed so the newest period list now reads
# in ascending order (2311, 2312, 2401, 2402, 2403, 2404, 2405) instead
# of the previous descending order (2405, 2404, 2403, 2402, 2401, 2312,
# 2311). The "Valor Mora" figures in column F for the first and last
# rows of that table (F16/F22) are swapped to match the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the "Periodo Mora" period labels (text values) in column E.
$ws.Range("E16").Value = "2311"
$ws.Range("E17").Value = "2312"
$ws.Range("E18").Value = "2401"
$ws.Range("E19").Value = "2402"
$ws.Range("E20").Value = "2403"
$ws.Range("E21").Value = "2404"
$ws.Range("E22").Value = "2405"

# Swap the "Valor Mora" amounts for the first/last rows to follow the
# new period ordering.
$ws.Range("F16").Value = 64000
$ws.Range("F22").Value = 46933
